$wb = $excel.ActiveWorkbook

# Add the new worksheet positioned right after the most recent existing week,
# so it lands at the end of the tab strip (sheetId 10, last position).
$srcHeader = $wb.Worksheets.Item("magapoke_2025-12-17")
$ws = $wb.Worksheets.Add($null, $srcHeader)
$ws.Name = "magapoke_2025-12-24"

# Copy header formatting/values (rank/title, bold+border+centered style) from prior week sheet
$srcHeader.Range("A1:B1").Copy($ws.Range("A1:B1"))

# Populate ranking rows
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "ブルーロック"
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "WIND BREAKER"
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！"
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "十字架のろくにん"
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "ガチアクタ"
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "東京卍リベンジャーズ"
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "みいちゃんと山田さん"
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "ベイビーステップ"
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "ギルティサークル"
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "島耕作"
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "FAIRY TAIL 100 YEARS QUEST"
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "薫る花は凛と咲く"
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "愛妻の裏アカ"
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "黄昏町プリズナーズ"
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "魔女と傭兵"
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "君が僕らを悪魔と呼んだ頃"
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "南海トラフ巨大地震"
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "転生したら第七王子だったので、気ままに魔術を極めます"
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "ハードワーカー中田"
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～"
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～"
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "蒼く染めろ"
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "K-9~警視庁公安部公安第9課異能対策係~"
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "アルキメデスの大戦"
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "異世界ウォーキング"
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す"
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "食糧人類-Starving Anonymous-"
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "グラぱらっ！"
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～"
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "せいぶつ部の田辺くん"
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "ハナバス　苔石花江のバスケ論"
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "幼馴染とはラブコメにならない"
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "ドラハチ"
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "我間乱 ―修羅―"
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "ナキナギ"
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "ひゃくえむ。"
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "ともだちづくり"
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "おやすみ ふみさん"
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "いじめるヤバイ奴"
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "さわらないで小手指くん"
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "イレギュラーズ"
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "アオバノバスケ"
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "ストーカー行為がバレて人生終了男"
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "屋根の下のアルテミス"
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！"
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜"
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "降り積もれ孤独な死よ"
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！"
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～"
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "なれの果ての僕ら"
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜"
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "黒猫と魔女の教室"
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "デスティニーラバーズ"
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "となりの黒川さん"
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。"
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "ジュミドロ"
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜"
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = "普通の本はありません！"
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "限界集落を脱村した錬金術士、都会で`"最強`"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～"
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "劣等人の魔剣使い　スキルボードを駆使して最強に至る"
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～"
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "時々ボソッとロシア語でデレる隣のアーリャさん"
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "インフェクション"
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～"
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "五輪の女神さま 〜なでしこ寮のメダルごはん〜"
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "金田一少年の事件簿外伝 犯人たちの事件簿"
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "剣帝学院の魔眼賢者"
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "阿武ノーマル"
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "MYS"
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "東京卍リベンジャーズ～場地圭介からの手紙～"
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "可愛いだけじゃない式守さん"
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = "春くらり"
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = "デッドアカウント"
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = "リスナーに騙されてダンジョンの最下層から脱出RTAすることになった"
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "DAYS外伝"
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～"
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = "中華一番！極"
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～"
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～"
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "私をセンターにすると誓いますか？"
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = "日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人"
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "田んぼで拾った女騎士、田舎で俺の嫁だと思われている"
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "君が監督！"
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = "この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～"
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = "ウイニング パス"
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～"
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = "人間消失"
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = "はじめの一歩"
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = "彼女、お借りします"
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = "勇者と呼ばれた後に　―そして無双男は家族を創る―"
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "七聖剣と魔剣の姫"
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = "東京ネオンスキャンダル"
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "FAIRY TAIL CITY HERO"
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "追放されなかった男　～二度目の人生は土下座から始まりました～"
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = "冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜"
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = "イジらないで、長瀞さん"
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～"
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = "Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜"
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = "堂島くんは動じない"
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = "捨てられた転生賢者～魔物の森で最強の大魔帝国を作り上げる～"

$ws.Range("A1").Select()

Write-Output "Added sheet magapoke_2025-12-24 with $($ws.UsedRange.Rows.Count) rows"
